$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.907.85"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "3.032.03"

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.96%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.525"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.48%  "

$ws.Range("D9").Value = "3.031.89"

$ws.Range("E10").Value = "  -2.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.75"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.444"
$ws.Range("D12").Style = "Normal"

$ws.Range("E13").Value = "  -1.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.92%  "

$ws.Range("E15").Value = "  +2.33%  "

$ws.Range("D16").Value = "3.532.02"
$ws.Range("E16").Value = "  -0.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").Value = "62.848.03"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").Value = "3.032.43"
$ws.Range("E19").Value = "  -0.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.693"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.47%  "

$ws.Range("E23").Value = "  -0.68%  "

$ws.Range("E24").Value = "  -1.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.72%  "

$ws.Range("E26").Value = "  -2.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.79%  "

$ws.Range("E28").Value = "  +0.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.49%  "

$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.97%  "

$ws.Range("E34").Value = "  -3.89%  "

$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("D36").Value = "0.0₃0804"
$ws.Range("E36").Value = "  -1.71%  "

$ws.Range("E37").Value = "  -2.46%  "

$ws.Range("E38").Value = "  -1.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "426.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.95%  "

$ws.Range("E43").Value = "  +2.79%  "

$ws.Range("E44").Value = "  -1.25%  "

$ws.Range("D45").Value = "2.805.35"
$ws.Range("E45").Value = "  +0.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0357"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("E47").Value = "  -8.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.78%  "

$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.80%  "

$ws.Range("E51").Value = "  -0.32%  "
